$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: D9 (title) and E9 (link)
$ws.Range("D9").Value = "학위 인증 후기 – 2.(좀 이상하지만) 괜찮은데 스위스?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/accreditation-procedure-experience-2/#utm_source=rss&utm_medium=rss&utm_campaign=accreditation-procedure-experience-2"

# Row 12: D12 (title) and E12 (link)
$ws.Range("D12").Value = "[파이썬 라이브러리를 활용한 머신러닝] 독자 리뷰를 소개합니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/05/16/%ed%8c%8c%ec%9d%b4%ec%8d%ac-%eb%9d%bc%ec%9d%b4%eb%b8%8c%eb%9f%ac%eb%a6%ac%eb%a5%bc-%ed%99%9c%ec%9a%a9%ed%95%9c-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-%eb%8f%85%ec%9e%90-%eb%a6%ac%eb%b7%b0%eb%a5%bc/"

# Row 26: D26 (title) only
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 37: D37 (title) and E37 (link)
$ws.Range("D37").Value = "[Paper Review] AugNLG: Few-shot Natural Language Generation using Self-trained Data Augmentation"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1981&mod=document&pageid=1"

# Row 51: D51 (title) and E51 (link)
$ws.Range("D51").Value = "[python+pandas] 여러 데이터프레임 하나의 엑셀 파일 내 각각 다른 시트에 저장하기"
$ws.Range("E51").Value = "https://bskyvision.com/1277"
